# Add 2022-Q3 data
# ------------------------------------------------------------------
# 1) Create the new "2022-Q3" sheet by duplicating the existing
#    "2022-Q2" sheet (so it inherits identical formatting/styles),
#    placing it immediately before "2022-Q2". Then trim/overwrite it
#    with the new quarter's fund-holding data.
# 2) Update the "总计" (summary) sheet: insert the new 2022-Q3 totals
#    at the top of the data list and push the rest of the list down
#    by one row, appending the row that drops off the end
#    (2020-Q4) as a brand new last row.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------- 1) New "2022-Q3" worksheet -----------------------------------

$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# The source sheet has 5 data rows (rows 2-6); the new sheet only needs 4
# (rows 2-5), so drop the extra trailing row.
$q3.Rows.Item(6).Delete()

# Row 1 (headers) is already correct/identical, carried over from the copy.

function Set-TextCell($sheet, $row, $col, $text) {
    $cell = $sheet.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# Row 2
$q3.Cells.Item(2, 1).Value = 0
Set-TextCell $q3 2 2 "012321"
$q3.Cells.Item(2, 3).Value = "东财中证云计算指数增强A"
Set-TextCell $q3 2 4 "1.31"
Set-TextCell $q3 2 5 "93.47"
Set-TextCell $q3 2 6 "3.82"
Set-TextCell $q3 2 7 "0.0500"
$q3.Cells.Item(2, 8).Value = 10

# Row 3
$q3.Cells.Item(3, 1).Value = 1
Set-TextCell $q3 3 2 "012322"
$q3.Cells.Item(3, 3).Value = "东财中证云计算指数增强C"
Set-TextCell $q3 3 4 "0.51"
Set-TextCell $q3 3 5 "93.47"
Set-TextCell $q3 3 6 "3.82"
Set-TextCell $q3 3 7 "0.0195"
$q3.Cells.Item(3, 8).Value = 10

# Row 4
$q3.Cells.Item(4, 1).Value = 2
Set-TextCell $q3 4 2 "006227"
$q3.Cells.Item(4, 3).Value = "华宝科技先锋混合A"
Set-TextCell $q3 4 4 "0.43"
Set-TextCell $q3 4 5 "90.67"
Set-TextCell $q3 4 6 "2.77"
Set-TextCell $q3 4 7 "0.0119"
$q3.Cells.Item(4, 8).Value = 10

# Row 5
$q3.Cells.Item(5, 1).Value = 3
Set-TextCell $q3 5 2 "010842"
$q3.Cells.Item(5, 3).Value = "华宝科技先锋混合C"
Set-TextCell $q3 5 4 "0.09"
Set-TextCell $q3 5 5 "90.67"
Set-TextCell $q3 5 6 "2.77"
Set-TextCell $q3 5 7 "0.0025"
$q3.Cells.Item(5, 8).Value = 10

# ---------- 2) "总计" summary sheet ---------------------------------------

$tot = $wb.Worksheets.Item("总计")

# Give new row 9 the same look (borders/font/number formats) as row 8
# before filling in its values, since it is a brand-new row.
$tot.Range("A8:D8").Copy()
$tot.Range("A9:D9").PasteSpecial(-4122)

# Column A is just a static 0-based counter per row and does not change.
# Columns B/C/D hold the quarterly figures; shift them all down by one
# row (oldest quarter falls into the newly-created row 9) and place the
# new 2022-Q3 figures at the top (row 2).

$tot.Cells.Item(2, 2).Value = "2022-Q3"
$tot.Cells.Item(2, 3).Value = 4
$tot.Cells.Item(2, 4).Value = 0.08

$tot.Cells.Item(3, 2).Value = "2022-Q2"
$tot.Cells.Item(3, 3).Value = 5
$tot.Cells.Item(3, 4).Value = 0.22

$tot.Cells.Item(4, 2).Value = "2022-Q1"
$tot.Cells.Item(4, 3).Value = 22
$tot.Cells.Item(4, 4).Value = 5.57

$tot.Cells.Item(5, 2).Value = "2021-Q4"
$tot.Cells.Item(5, 3).Value = 16
$tot.Cells.Item(5, 4).Value = 2.67

$tot.Cells.Item(6, 2).Value = "2021-Q3"
$tot.Cells.Item(6, 3).Value = 12
$tot.Cells.Item(6, 4).Value = 6.1

$tot.Cells.Item(7, 2).Value = "2021-Q2"
$tot.Cells.Item(7, 3).Value = 23
$tot.Cells.Item(7, 4).Value = 9.24

$tot.Cells.Item(8, 2).Value = "2021-Q1"
$tot.Cells.Item(8, 3).Value = 17
$tot.Cells.Item(8, 4).Value = 9.16

$tot.Cells.Item(9, 1).Value = 7
$tot.Cells.Item(9, 2).Value = "2020-Q4"
$tot.Cells.Item(9, 3).Value = 30
$tot.Cells.Item(9, 4).Value = 13.47
